$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Create a "Bad" cell style (red text on light-red fill) so it shows up in
# the workbook's style gallery alongside the existing "Good" style, and so
# the cells below can reference it by name.
# ---------------------------------------------------------------------------
$badStyle = $wb.Styles.Add("Bad")
$badStyle.Font.Size = 12
$badStyle.Font.Color = 393372        # BGR for RGB FF9C0006
$badStyle.Interior.Color = 13551615  # BGR for RGB FFFFC7CE

# ---------------------------------------------------------------------------
# Existing rows 9 & 10: the expected result changed from "passed" to "fail",
# so the cell style flips from "Good" to the new "Bad" style.
# ---------------------------------------------------------------------------
$ws.Range("C9").Value = "fail"
$ws.Range("C9").Style = "Bad"
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108

$ws.Range("C10").Value = "fail"
$ws.Range("C10").Style = "Bad"
$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("C10").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# New test-case rows 11-15.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = 1.9
$ws.Range("B11").Value = 'Попытка создания проекта при заполнении поля "Описание" невалидным значением'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "fail"
$ws.Range("C11").Style = "Bad"
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("C11").VerticalAlignment = -4108

$ws.Range("A12").Value = 1.11
$ws.Range("B12").Value = "Попытка создания проекта с заполнением поля название невалидными символами"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "fail"
$ws.Range("C12").Style = "Bad"
$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("C12").VerticalAlignment = -4108

$ws.Range("A13").Value = 1.12
$ws.Range("B13").Value = "Попытка создания проекта с заполнением поля Сокращенное название невалидными символами"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "fail"
$ws.Range("C13").Style = "Bad"
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("C13").VerticalAlignment = -4108

$ws.Range("A14").Value = 1.13
$ws.Range("B14").Value = "Попытка создания проекта с заполнением поля описание невалидными символами"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "fail"
$ws.Range("C14").Style = "Bad"
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("C14").VerticalAlignment = -4108

$ws.Range("A15").Value = 1.14
$ws.Range("B15").Value = 'Попытка создания проекта при незаполненом поле "Описание проекта"'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "passed"
$ws.Range("C15").Style = "Good"
$ws.Range("C15").HorizontalAlignment = -4108
$ws.Range("C15").VerticalAlignment = -4108

# Re-apply the wrap/valign formatting (style = "Normal" above resets number
# format / alignment) used by the rest of column A & B.
$ws.Range("A11:A15").NumberFormat = "00000"
$ws.Range("A3:A15").WrapText = $true
$ws.Range("A3:A15").VerticalAlignment = -4160
$ws.Range("B3:B15").WrapText = $true
$ws.Range("B3:B15").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Row heights (pt) - mirrors the auto height Excel computed for the wrapped
# text once the new rows were inserted.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 34
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 34
$ws.Rows.Item(10).RowHeight = 34
$ws.Rows.Item(11).RowHeight = 34
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(13).RowHeight = 34
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 34

# ---------------------------------------------------------------------------
# Selection, matching the cursor position left behind in the saved file.
# ---------------------------------------------------------------------------
$ws.Range("F13").Select()

Write-Host "done"
